$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2236.25
$ws.Range("I62").Value = 1978
$ws.Range("K62").Value = 1978
$ws.Range("M62").Value = -1354

$ws.Range("H65").Value = 2236.25
$ws.Range("I65").Value = 1978
$ws.Range("K65").Value = 9890
$ws.Range("M65").Value = -6770

$ws.Range("H132").Value = 4235.6284
$ws.Range("I132").Value = 4158.2666
$ws.Range("J132").Value = 4699.8
$ws.Range("K132").Value = 12474.7998
$ws.Range("L132").Value = 14099.4
$ws.Range("M132").Value = -9944.799800000001
$ws.Range("N132").Value = -19159.4

$ws.Range("H135").Value = 699.7646999999999
$ws.Range("I135").Value = 619.4815
$ws.Range("J135").Value = 1009.4286
$ws.Range("K135").Value = 5575.3335
$ws.Range("L135").Value = 9084.857399999999
$ws.Range("M135").Value = -3040.3335
$ws.Range("N135").Value = -14154.8574

$ws.Range("H137").Value = 905.48
$ws.Range("I137").Value = 746.6896400000001
$ws.Range("J137").Value = 1124.762
$ws.Range("K137").Value = 2240.06892
$ws.Range("L137").Value = 3374.286
$ws.Range("M137").Value = 309.9310799999998
$ws.Range("N137").Value = -8474.286

$ws.Range("H138").Value = 1249.5
$ws.Range("I138").Value = 539.8043
$ws.Range("J138").Value = 1854.0555
$ws.Range("K138").Value = 1619.4129
$ws.Range("L138").Value = 5562.166499999999
$ws.Range("M138").Value = 3520.5871
$ws.Range("N138").Value = -15842.1665

$ws.Range("H141").Value = 2778.1667
$ws.Range("I141").Value = 934.9286
$ws.Range("K141").Value = 2804.7858
$ws.Range("M141").Value = 2375.2142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2232.97
$ws.Range("I32").Value = 2049.811
$ws.Range("J32").Value = 3881.4
$ws.Range("K32").Value = 2049.811
$ws.Range("L32").Value = 3881.4
$ws.Range("M32").Value = -1762.811
$ws.Range("N32").Value = -4455.4

$ws.Range("H74").Value = 2044.0238
$ws.Range("I74").Value = 1199.7222
$ws.Range("J74").Value = 2677.25
$ws.Range("K74").Value = 1199.7222
$ws.Range("L74").Value = 2677.25
$ws.Range("M74").Value = -325.7221999999999
$ws.Range("N74").Value = -4425.25

$ws.Range("H77").Value = 2044.0238
$ws.Range("I77").Value = 1199.7222
$ws.Range("J77").Value = 2677.25
$ws.Range("K77").Value = 5998.611
$ws.Range("L77").Value = 13386.25
$ws.Range("M77").Value = -1630.611
$ws.Range("N77").Value = -22122.25

$ws.Range("H102").Value = 1673.3334
$ws.Range("I102").Value = 1673.3334
$ws.Range("K102").Value = 1673.3334
$ws.Range("M102").Value = -51.33339999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3690.9
$ws.Range("I20").Value = 3460
$ws.Range("J20").Value = 3921.8
$ws.Range("K20").Value = 3460
$ws.Range("L20").Value = 3921.8
$ws.Range("M20").Value = -3213
$ws.Range("N20").Value = -4415.8

$ws.Range("H99").Value = 1416.6666
$ws.Range("I99").Value = 1125
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 1125
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = 373
$ws.Range("N99").Value = -4996

$ws.Range("H134").Value = 2992.2727
$ws.Range("I134").Value = 2944.2856
$ws.Range("K134").Value = 8832.856800000001
$ws.Range("M134").Value = -6297.856800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4662.3
$ws.Range("I31").Value = 1525.8695
$ws.Range("J31").Value = 6611.973
$ws.Range("K31").Value = 1525.8695
$ws.Range("L31").Value = 6611.973
$ws.Range("M31").Value = -1230.8695
$ws.Range("N31").Value = -7201.973

$ws.Range("H34").Value = 4662.3
$ws.Range("I34").Value = 1525.8695
$ws.Range("J34").Value = 6611.973
$ws.Range("K34").Value = 1525.8695
$ws.Range("L34").Value = 6611.973
$ws.Range("M34").Value = -1323.8695
$ws.Range("N34").Value = -7015.973

$ws.Range("H132").Value = 6174943.5
$ws.Range("I132").Value = 1877.1177
$ws.Range("J132").Value = 16669157
$ws.Range("K132").Value = 5631.3531
$ws.Range("L132").Value = 50007471
$ws.Range("M132").Value = -3101.3531
$ws.Range("N132").Value = -50012531

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1402.2903
$ws.Range("I5").Value = 457.11765
$ws.Range("J5").Value = 2550
$ws.Range("K5").Value = 1371.35295
$ws.Range("L5").Value = 7650
$ws.Range("M5").Value = -1259.35295
$ws.Range("N5").Value = -7874

$ws.Range("H121").Value = 514034.62
$ws.Range("I121").Value = 5000485
$ws.Range("J121").Value = 1297.4286
$ws.Range("K121").Value = 15001455
$ws.Range("L121").Value = 3892.2858
$ws.Range("M121").Value = -15000145
$ws.Range("N121").Value = -6512.2858

$ws.Range("H122").Value = 2565.5881
$ws.Range("I122").Value = 386.51724
$ws.Range("J122").Value = 5438
$ws.Range("K122").Value = 3478.65516
$ws.Range("L122").Value = 48942
$ws.Range("M122").Value = -1028.65516
$ws.Range("N122").Value = -53842

$ws.Range("H127").Value = 9615735
$ws.Range("J127").Value = 9615735
$ws.Range("L127").Value = 28847205
$ws.Range("N127").Value = -28857125

$ws.Range("H135").Value = 1402.2903
$ws.Range("I135").Value = 457.11765
$ws.Range("J135").Value = 2550
$ws.Range("K135").Value = 4114.05885
$ws.Range("L135").Value = 22950
$ws.Range("M135").Value = -1579.05885
$ws.Range("N135").Value = -28020

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 6829.5713
$ws.Range("I29").Value = 507
$ws.Range("J29").Value = 7883.3335
$ws.Range("K29").Value = 507
$ws.Range("L29").Value = 7883.3335
$ws.Range("M29").Value = -217
$ws.Range("N29").Value = -8463.333500000001

$ws.Range("H97").Value = 568.3333
$ws.Range("I97").Value = 579.61536
$ws.Range("J97").Value = 495
$ws.Range("K97").Value = 579.61536
$ws.Range("L97").Value = 495
$ws.Range("M97").Value = -83.61536000000001
$ws.Range("N97").Value = -1487

$ws.Range("H122").Value = 2599.8
$ws.Range("I122").Value = 2333
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 6999
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -4549
$ws.Range("N122").Value = -13900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3379.3428
$ws.Range("I22").Value = 432.46155
$ws.Range("J22").Value = 5120.6816
$ws.Range("K22").Value = 432.46155
$ws.Range("L22").Value = 5120.6816
$ws.Range("M22").Value = -137.46155
$ws.Range("N22").Value = -5710.6816

$ws.Range("H27").Value = 3379.3428
$ws.Range("I27").Value = 432.46155
$ws.Range("J27").Value = 5120.6816
$ws.Range("K27").Value = 432.46155
$ws.Range("L27").Value = 5120.6816
$ws.Range("M27").Value = -325.46155
$ws.Range("N27").Value = -5334.6816

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4488676.5
$ws.Range("I132").Value = 1569.6842
$ws.Range("J132").Value = 10803864
$ws.Range("K132").Value = 4709.0526
$ws.Range("L132").Value = 32411592
$ws.Range("M132").Value = -2179.0526
$ws.Range("N132").Value = -32416652
